$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values that could be misinterpreted as numbers by Excel
# are protected by temporarily forcing text format, then the style is
# reset back to Normal so no stray formatting is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.254.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.355.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.650"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.56"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0978"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "27.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.865"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.361.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.244.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "249.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.128"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0691"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0958"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.441.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.578.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000202"
$ws.Range("D51").Style = "Normal"

# Remaining plain-text updates (names, links, padded percentage strings)
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("E3").Value = "  +5.28%  "
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("E7").Value = "  +13.54%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +11.90%  "
$ws.Range("E10").Value = "  +2.06%  "
$ws.Range("E11").Value = "  +3.45%  "
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("E13").Value = "  +5.18%  "
$ws.Range("E14").Value = "  +4.76%  "
$ws.Range("E15").Value = "  +5.03%  "
$ws.Range("E16").Value = "  +4.78%  "
$ws.Range("E17").Value = "  +4.37%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("E19").Value = "  +5.10%  "
$ws.Range("E20").Value = "  +4.91%  "
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  +2.69%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E26").Value = "  +3.66%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("E28").Value = "  +4.48%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("E30").Value = "  +8.27%  "
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("E33").Value = "  +3.40%  "
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("E35").Value = "  +4.69%  "
$ws.Range("E36").Value = "  +3.34%  "
$ws.Range("E37").Value = "  +4.24%  "
$ws.Range("E38").Value = "  +7.66%  "
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E40").Value = "  +12.19%  "
$ws.Range("B41").Value = "BinanceUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +4.76%  "
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("E44").Value = "  +9.48%  "
$ws.Range("E45").Value = "  +2.16%  "
$ws.Range("E46").Value = "  +2.77%  "
$ws.Range("E47").Value = "  +2.91%  "
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("E49").Value = "  +5.36%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("E51").Value = "  -1.97%  "
